# Auto-generated edit script: apply numeric cell updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 429.35715
$ws.Range("I2").Value = 291.5
$ws.Range("K2").Value = 291.5
$ws.Range("M2").Value = -178.5

# Row 15
$ws.Range("H15").Value = 781.8022
$ws.Range("I15").Value = 781.8022
$ws.Range("K15").Value = 2345.4066
$ws.Range("M15").Value = -2176.4066

# Row 32
$ws.Range("H32").Value = 5538.091
$ws.Range("J32").Value = 5927.375
$ws.Range("L32").Value = 5927.375
$ws.Range("N32").Value = -6579.375

# Row 86
$ws.Range("H86").Value = 11809.25
$ws.Range("I86").Value = 10995.667
$ws.Range("J86").Value = 14250
$ws.Range("K86").Value = 10995.667
$ws.Range("L86").Value = 14250
$ws.Range("M86").Value = -9872.666999999999
$ws.Range("N86").Value = -16496

# Row 89
$ws.Range("H89").Value = 11809.25
$ws.Range("I89").Value = 10995.667
$ws.Range("J89").Value = 14250
$ws.Range("K89").Value = 54978.335
$ws.Range("L89").Value = 71250
$ws.Range("M89").Value = -49362.335
$ws.Range("N89").Value = -82482

# Row 112
$ws.Range("H112").Value = 5749910.5
$ws.Range("J112").Value = 6026954.5
$ws.Range("L112").Value = 18080863.5
$ws.Range("N112").Value = -18083079.5

# Row 131
$ws.Range("H131").Value = 1197.1111
$ws.Range("I131").Value = 1182.1428
$ws.Range("K131").Value = 3546.4284
$ws.Range("M131").Value = 1493.5716

# Row 138
$ws.Range("H138").Value = 4336443
$ws.Range("J138").Value = 4702478
$ws.Range("L138").Value = 14107434
$ws.Range("N138").Value = -14117714

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10995091
$ws.Range("I32").Value = 11368821
$ws.Range("J32").Value = 32332.666
$ws.Range("K32").Value = 11368821
$ws.Range("L32").Value = 32332.666
$ws.Range("M32").Value = -11368534
$ws.Range("N32").Value = -32906.666

# Row 110
$ws.Range("H110").Value = 24342.234
$ws.Range("I110").Value = 25551.125
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 25551.125
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -23506.125
$ws.Range("N110").Value = -9090

# Row 132
$ws.Range("H132").Value = 34486948
$ws.Range("I132").Value = 4207.7036
$ws.Range("K132").Value = 12623.1108
$ws.Range("M132").Value = -10093.1108

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2763.375
$ws.Range("I20").Value = 3592.3333
$ws.Range("J20").Value = 1697.5714
$ws.Range("K20").Value = 3592.3333
$ws.Range("L20").Value = 1697.5714
$ws.Range("M20").Value = -3345.3333
$ws.Range("N20").Value = -2191.5714

# Row 86
$ws.Range("H86").Value = 34439.1
$ws.Range("I86").Value = 52203
$ws.Range("J86").Value = 29998.125
$ws.Range("K86").Value = 52203
$ws.Range("L86").Value = 29998.125
$ws.Range("M86").Value = -51080
$ws.Range("N86").Value = -32244.125

# Row 89
$ws.Range("H89").Value = 34439.1
$ws.Range("I89").Value = 52203
$ws.Range("J89").Value = 29998.125
$ws.Range("K89").Value = 261015
$ws.Range("L89").Value = 149990.625
$ws.Range("M89").Value = -255399
$ws.Range("N89").Value = -161222.625

# Row 99
$ws.Range("H99").Value = 2642.5652
$ws.Range("I99").Value = 1351.2
$ws.Range("J99").Value = 5063.875
$ws.Range("K99").Value = 1351.2
$ws.Range("L99").Value = 5063.875
$ws.Range("M99").Value = 146.8
$ws.Range("N99").Value = -8059.875

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 6133
$ws.Range("J62").Value = 6133
$ws.Range("L62").Value = 6133
$ws.Range("N62").Value = -7381

# Row 65
$ws.Range("H65").Value = 6133
$ws.Range("J65").Value = 6133
$ws.Range("L65").Value = 30665
$ws.Range("N65").Value = -36905

# Row 70
$ws.Range("H70").Value = 18000
$ws.Range("J70").Value = 18000
$ws.Range("L70").Value = 18000
$ws.Range("N70").Value = -18630

# Row 73
$ws.Range("H73").Value = 18000
$ws.Range("J73").Value = 18000
$ws.Range("L73").Value = 18000
$ws.Range("N73").Value = -20184

# Row 99
$ws.Range("H99").Value = 7457
$ws.Range("I99").Value = 7828.0835
$ws.Range("J99").Value = 5230.5
$ws.Range("K99").Value = 7828.0835
$ws.Range("L99").Value = 5230.5
$ws.Range("M99").Value = -6330.0835
$ws.Range("N99").Value = -8226.5

# Row 126
$ws.Range("H126").Value = 7457
$ws.Range("I126").Value = 7828.0835
$ws.Range("J126").Value = 5230.5
$ws.Range("K126").Value = 23484.2505
$ws.Range("L126").Value = 15691.5
$ws.Range("M126").Value = -21014.2505
$ws.Range("N126").Value = -20631.5

# Row 132
$ws.Range("H132").Value = 3088.353
$ws.Range("I132").Value = 3163.25
$ws.Range("J132").Value = 1890
$ws.Range("K132").Value = 9489.75
$ws.Range("L132").Value = 5670
$ws.Range("M132").Value = -6959.75
$ws.Range("N132").Value = -10730

# Row 134
$ws.Range("H134").Value = 1652.2106
$ws.Range("I134").Value = 1523.6875
$ws.Range("J134").Value = 2337.6667
$ws.Range("K134").Value = 4571.0625
$ws.Range("L134").Value = 7013.000100000001
$ws.Range("M134").Value = -2036.0625
$ws.Range("N134").Value = -12083.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 67325.734
$ws.Range("I2").Value = 1416.6666
$ws.Range("J2").Value = 111265.11
$ws.Range("K2").Value = 8499.999599999999
$ws.Range("L2").Value = 667590.66
$ws.Range("M2").Value = -8386.999599999999
$ws.Range("N2").Value = -667816.66

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 113
$ws.Range("H113").Value = 4792.1333
$ws.Range("J113").Value = 5788.3
$ws.Range("L113").Value = 17364.9
$ws.Range("N113").Value = -21704.9

# Row 140
$ws.Range("H140").Value = 2490.3845
$ws.Range("I140").Value = 2365
$ws.Range("K140").Value = 7095
$ws.Range("M140").Value = -1915

# Row 141
$ws.Range("H141").Value = 7071.7036
$ws.Range("I141").Value = 4864
$ws.Range("J141").Value = 10282.909
$ws.Range("K141").Value = 14592
$ws.Range("L141").Value = 30848.727
$ws.Range("M141").Value = -9412
$ws.Range("N141").Value = -41208.727

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5370.125
$ws.Range("I80").Value = 5086.9165
$ws.Range("K80").Value = 5086.9165
$ws.Range("M80").Value = -4088.9165

# Row 83
$ws.Range("H83").Value = 5370.125
$ws.Range("I83").Value = 5086.9165
$ws.Range("K83").Value = 25434.5825
$ws.Range("M83").Value = -20442.5825

# Row 132
$ws.Range("H132").Value = 4214.825
$ws.Range("I132").Value = 3803.0278
$ws.Range("K132").Value = 11409.0834
$ws.Range("M132").Value = -8879.0834

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1619.95

# Row 122
$ws.Range("H122").Value = 5970.3335
$ws.Range("I122").Value = 5115.727
$ws.Range("J122").Value = 7313.2856
$ws.Range("K122").Value = 15347.181
$ws.Range("L122").Value = 21939.8568
$ws.Range("M122").Value = -12897.181
$ws.Range("N122").Value = -26839.8568

# Row 131
$ws.Range("H131").Value = 54127
$ws.Range("J131").Value = 54127
$ws.Range("L131").Value = 54127
$ws.Range("N131").Value = -64207

# Row 136
$ws.Range("H136").Value = 8099.3335
$ws.Range("I136").Value = 7149
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 21447
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -18897
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 111226110
$ws.Range("I122").Value = 143004180
$ws.Range("J122").Value = 2887.5
$ws.Range("K122").Value = 429012540
$ws.Range("L122").Value = 8662.5
$ws.Range("M122").Value = -429010090
$ws.Range("N122").Value = -13562.5

# Row 126
$ws.Range("H126").Value = 6634.905
$ws.Range("I126").Value = 7037.706
$ws.Range("J126").Value = 4923
$ws.Range("K126").Value = 21113.118
$ws.Range("L126").Value = 14769
$ws.Range("M126").Value = -18643.118
$ws.Range("N126").Value = -19709

# Row 136
$ws.Range("H136").Value = 1384.9395
$ws.Range("I136").Value = 1334.7693
$ws.Range("K136").Value = 4004.3079
$ws.Range("M136").Value = -1454.3079
